$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 11 and 12: replace the SUM(...) formulas in column A with plain values
$ws.Range("A11").Value = 120000
$ws.Range("A12").Value = 120000

# Move the "I alt" total row from row 15 up to row 13 (rows 13-14 were empty)
$ws.Range("D13").Formula = $ws.Range("D15").Formula
$ws.Range("E13").Value = $ws.Range("E15").Value2
$ws.Range("D15").ClearContents()
$ws.Range("E15").ClearContents()

# Update the visible selection/scroll position to match the edited area
$ws.Range("A10:E13").Select()
